$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# 1) "2016-08-31 06:58:20" -> "2016-08-31 07:00:38"
#    Overview!G2, Overview!G3, de-de!H2, de-de!H3
$wsOverview.Range("G2").Value = "2016-08-31 07:00:38"
$wsOverview.Range("G3").Value = "2016-08-31 07:00:38"
$wsDeDe.Range("H2").Value = "2016-08-31 07:00:38"
$wsDeDe.Range("H3").Value = "2016-08-31 07:00:38"

# 2) "ht" -> "mt"
#    zh-cn!E2, zh-cn!E3, de-de!E2, de-de!E3
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# 3) "2016-08-31 06:57:59" -> "2016-08-31 07:00:13"
#    zh-cn!H2, zh-cn!H3
$wsZhCn.Range("H2").Value = "2016-08-31 07:00:13"
$wsZhCn.Range("H3").Value = "2016-08-31 07:00:13"

# 4) "2016-08-31 06:59:00" -> "2016-08-31 07:01:41"
#    zh-cn!K2, zh-cn!K3
$wsZhCn.Range("K2").Value = "2016-08-31 07:01:41"
$wsZhCn.Range("K3").Value = "2016-08-31 07:01:41"

# 5) "2016-08-31 06:59:24" -> "2016-08-31 07:02:00"
#    de-de!K2, de-de!K3
$wsDeDe.Range("K2").Value = "2016-08-31 07:02:00"
$wsDeDe.Range("K3").Value = "2016-08-31 07:02:00"
